$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert two new columns (firstCond / secondCond) before the existing
# firstSkill/secondSkill columns in the left data block (B:C), and two more
# new columns for the corresponding lookup formulas in the right-hand
# formula block (K:L, which sits right after the shifted subNr-concat
# column).
# ---------------------------------------------------------------------------
$ws.Columns("B:C").Insert()
$ws.Columns("K:L").Insert()

# ---------------------------------------------------------------------------
# Header row (row 1)
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "firstCond"
$ws.Range("F1").Value = "firstScore"
$ws.Range("G1").Value = "secondScore"

$ws.Range("K1").Value = "firstCond"

$ws.Range("B2").Value = "Teaching"
$ws.Range("C2").Value = "Performing"

$ws.Range("C1").Value = "secondCond"
$ws.Range("L1").Value = "secondCond"

# ---------------------------------------------------------------------------
# Data rows 2-9: firstCond / secondCond alternate Teaching / Performing
# ---------------------------------------------------------------------------
$firstConds = @("Teaching","Performing","Teaching","Performing","Teaching","Performing","Teaching","Performing")
$secondConds = @("Performing","Teaching","Performing","Teaching","Performing","Teaching","Performing","Teaching")

for ($i = 0; $i -lt 8; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $firstConds[$i]
    $ws.Cells.Item($row, 3).Value = $secondConds[$i]
}

# ---------------------------------------------------------------------------
# Lookup-formula columns K (=firstCond) / L (=secondCond) for rows 2-9.
# Row 2 gets its own (non-shared) formula; rows 3-9 share one formula group,
# matching the pattern already present for the other lookup columns.
# ---------------------------------------------------------------------------
$ws.Range("K2").Formula = "=B2"
$ws.Range("L2").Formula = "=C2"
$ws.Range("K3:K9").Formula = "=B3"
$ws.Range("L3:L9").Formula = "=C3"

# The two newly inserted lookup columns broke up the shared-formula groups
# that used to cover the old H:L block (now M:P) - rebuild them the same way.
$ws.Range("M3:M9").Formula = "=D3"
$ws.Range("N3:N9").Formula = "=E3"
$ws.Range("O3:O9").Formula = "=F3"
$ws.Range("P3:P9").Formula = "=CONCATENATE(G3,I3)"

# ---------------------------------------------------------------------------
# Selection, to match the committed file's sheet view
# ---------------------------------------------------------------------------
$ws.Range("J15").Select()
